$wb = $excel.ActiveWorkbook

# Update the "About" sheet text: "For the U.S.:" -> "For the EU.:"
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A10").Value = "For the EU.:"

# Leave the cursor on A11 in the About sheet (next to the edited cell),
# then move to / leave active the MOU-large sheet.
[void]$wsAbout.Range("A11").Select()

$wsLarge = $wb.Worksheets.Item("MOU-large")
[void]$wsLarge.Activate()
